# Actualizacion automatica 2025-06-05 10:49:05
# Adds the "CUMPLIMIENTO MENSUAL" sheet (monthly compliance / goal-attainment
# report) built from the PRESUPUESTO (budget) vs VENTA (sales) figures that
# already live on "VENTAS POR GRUPO".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Create the new worksheet as the LAST tab in the workbook.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# Match the outline / summary settings used on the rest of the workbook.
$ws.Outline.SummaryRow    = 1
$ws.Outline.SummaryColumn = 1

# Page margins identical to the other sheets (0.75" L/R, 1" T/B, 0.5" header/footer).
$ws.PageSetup.LeftMargin   = 54
$ws.PageSetup.RightMargin  = 54
$ws.PageSetup.TopMargin    = 72
$ws.PageSetup.BottomMargin = 72
$ws.PageSetup.HeaderMargin = 36
$ws.PageSetup.FooterMargin = 36

# ---------------------------------------------------------------------
# 2. Column widths.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 31.16666666666667  # A -> 32
$ws.Columns.Item(2).ColumnWidth = 21.16666666666667  # B -> 22
$ws.Columns.Item(3).ColumnWidth = 21.16666666666667  # C -> 22
$ws.Columns.Item(4).ColumnWidth = 12.16666666666667  # D -> 13
$ws.Columns.Item(5).ColumnWidth = 21.16666666666667  # E -> 22
$ws.Columns.Item(6).ColumnWidth = 24.16666666666667  # F -> 25

# ---------------------------------------------------------------------
# 3. Header row.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

# Re-use the same bold / bordered / centered header style already used on
# the other two sheets instead of inventing a new one.
$ws1.Range("A1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Detail rows: one per product GRUPO for this ASESOR.
#    PRESUPUESTO (budget) and VENTA (sales-to-date) are the source figures;
#    POR CUMPLIR = PRESUPUESTO - VENTA, CUMPLIMIENTO = VENTA / PRESUPUESTO.
# ---------------------------------------------------------------------
$asesor = "GUERRERO FAREZ FABIAN MAURICIO"

$grupos = @(
    @("240X120 PORCELANATO", 9970.34304517915, 0),
    @("240X80 PORCELANATO", 27457.0076, 0),
    @("FREGADEROS DE COCINA", 1003, 0),
    @("GRANITO", 238.32, 0),
    @("GRIFERIAS", 106.82, 0),
    @("INODOROS", 1400, 0),
    @("LAVABOS", 1000, 0),
    @("LED", 300, 0),
    @("NO RESURTIBLES", 1300.5, 0),
    @("OTROS", 0, 0),
    @("PANELES DECORATIVOS", 350, 0),
    @("PANELES PU", 230, 0),
    @("PANELES PVC", 966, 0),
    @("PIEDRA SINTERIZADA", 13500, 0),
    @("PORCELANATO", 32741.45, 2529.84),
    @("PUERTAS DE SEGURIDAD", 684, 0),
    @("SAL SOLUBLE", 3200, 0)
)

$row = 2
$totalPresupuesto = 0
$totalVenta = 0

foreach ($g in $grupos) {
    $grupo = $g[0]
    $presupuesto = $g[1]
    $venta = $g[2]
    $porCumplir = $presupuesto - $venta
    if ($presupuesto -ne 0) {
        $cumplimiento = $venta / $presupuesto
    } else {
        $cumplimiento = 0
    }

    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $grupo
    $ws.Cells.Item($row, 3).Value = $presupuesto
    $ws.Cells.Item($row, 4).Value = $venta
    $ws.Cells.Item($row, 5).Value = $porCumplir
    $ws.Cells.Item($row, 6).Value = $cumplimiento

    $totalPresupuesto = $totalPresupuesto + $presupuesto
    $totalVenta = $totalVenta + $venta

    $row = $row + 1
}

$lastDataRow = $row - 1
$totalRow = $row

# ---------------------------------------------------------------------
# 5. TOTAL row.
# ---------------------------------------------------------------------
$totalPorCumplir = $totalPresupuesto - $totalVenta
$totalCumplimiento = $totalVenta / $totalPresupuesto

$ws.Cells.Item($totalRow, 2).Value = "TOTAL"
$ws.Cells.Item($totalRow, 3).Value = $totalPresupuesto
$ws.Cells.Item($totalRow, 4).Value = $totalVenta
$ws.Cells.Item($totalRow, 5).Value = $totalPorCumplir
$ws.Cells.Item($totalRow, 6).Value = $totalCumplimiento

# ---------------------------------------------------------------------
# 6. Number formats.
#    C/D/E -> currency ("$"#,##0.00, same custom format already used
#    elsewhere in the workbook); F -> percentage (0.00%).
# ---------------------------------------------------------------------
$ws.Range("C2:E" + $totalRow).NumberFormat = """$""#,##0.00"
$ws.Range("F2:F" + $totalRow).NumberFormat = "0.00%"

# Right-align the "TOTAL" label.
$ws.Cells.Item($totalRow, 2).HorizontalAlignment = -4152  # xlRight

Write-Host "CUMPLIMIENTO MENSUAL sheet created with" ($lastDataRow - 1) "group rows plus header and total."
